# Fill in the gap for the week of 43157 (row 7): Tuesday (column E) hours = 2.25.
# This is a missing entry that the shared SUM formula in column I (and the
# grand total in I19) will automatically pick up once recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 2.25

# Move the active selection to F9, matching where the author left off editing.
$ws.Range("F9").Select()

$excel.Calculate()
